# Appends a new weekly "03-03" A/0 column pair (AR/AS) to the sheet, carried
# forward from the previous "03-02" pair (AP/AQ), and normalizes the AQ
# column's stored type from inline text to a real number along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new date-pair labels in AR1 / AS1, matching AQ1's style ---
$ws.Range("AQ1").Copy($ws.Range("AR1"))
$ws.Range("AQ1").Copy($ws.Range("AS1"))
$ws.Range("AR1").Value = "03-03_A"
$ws.Range("AS1").Value = "03-03_0"

# --- Data rows: duplicate the AP/AQ pair forward into AR/AS ---
for ($r = 2; $r -le 173; $r++) {
    $apAddr = "AP" + $r
    $aqAddr = "AQ" + $r
    $arAddr = "AR" + $r
    $asAddr = "AS" + $r

    # AR gets AP's value + style (single Copy call carries both).
    $ws.Range($apAddr).Copy($ws.Range($arAddr))

    # AS gets AQ's (pre-conversion) value + style.
    $ws.Range($aqAddr).Copy($ws.Range($asAddr))

    # Normalize AQ from inline-string to a real number, when it actually
    # holds a value (truly-blank rows keep their empty inline-string cell).
    $aqVal = $ws.Range($aqAddr).Value()
    if (-not ($aqVal -eq $null -or $aqVal -eq "")) {
        $ws.Range($aqAddr).Value = [double]$aqVal
    }
}
